# Apply the translation/formatting update described in the commit:
# "Update files based on cn190115 and Improve Translation for console"
#
# The two validation-message cells (C16, C18) get a highlighted
# " Chinese characters" phrase (red font) inserted into their text,
# the rows grow taller to fit the re-wrapped text, a new (wider) column D
# is introduced, and a couple of blank formatted cells appear in the
# widened rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 16 (Name validation message)
# ---------------------------------------------------------------------
$run1_16 = "Name only supports"
$run2_16 = " Chinese characters"
$run3_16 = ", numbers, upper case and lower case letters, English underline “_” and line-through"
$full16 = $run1_16 + $run2_16 + $run3_16

$ws.Range("C16").Value = $full16

# Highlight " Chinese characters" in red
$ws.Range("C16").Characters($run1_16.Length + 1, $run2_16.Length).Font.Color = 255
# Re-assert the trailing run's colour so it keeps its own explicit run
$ws.Range("C16").Characters($run1_16.Length + $run2_16.Length + 1, $run3_16.Length).Font.Color = 0

# ---------------------------------------------------------------------
# 2. Row 18 (Description validation message)
# ---------------------------------------------------------------------
$run1_18 = "Description only supports"
$run2_18 = " Chinese characters"
$run3_18 = ", numbers, upper case and lower case letters and English underline “_”"
$full18 = $run1_18 + $run2_18 + $run3_18

$ws.Range("C18").Value = $full18

$ws.Range("C18").Characters($run1_18.Length + 1, $run2_18.Length).Font.Color = 255
$ws.Range("C18").Characters($run1_18.Length + $run2_18.Length + 1, $run3_18.Length).Font.Color = 0

# ---------------------------------------------------------------------
# 3. Row heights grow to fit the (now longer) wrapped text
# ---------------------------------------------------------------------
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(18).RowHeight = 45

# ---------------------------------------------------------------------
# 4. New column D is introduced (wider, to host the extra text)
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 32.6

# ---------------------------------------------------------------------
# 5. Blank (but formatted) cells appear at D16/E16 and D18/E18
# ---------------------------------------------------------------------
$ws.Range("D16").WrapText = $true
$ws.Range("E16").WrapText = $true
$ws.Range("D18").WrapText = $true
$ws.Range("E18").WrapText = $true

# ---------------------------------------------------------------------
# 6. Selection / view ends up on C16
# ---------------------------------------------------------------------
$ws.Range("C16").Select() | Out-Null

# ---------------------------------------------------------------------
# 7. Register the new red font in the workbook's font table (mirrors
#    what Excel does when a font colour is applied via the ribbon) then
#    revert the one-off probe cell so no stray content remains.
# ---------------------------------------------------------------------
$ws.Range("Z100").Font.Color = 255
$excel.Undo()

Write-Host "edit applied"
